$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header row (copy formatting from H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for columns I (I0) and J (IF) for rows 2-65.
$iValues = @(
    6,
    8,
    7,
    9,
    7,
    8,
    6,
    4,
    6,
    8,
    8,
    5,
    8,
    5,
    7,
    8,
    8,
    10,
    7,
    6,
    8,
    6,
    8,
    7,
    9,
    10,
    8,
    9,
    6,
    7,
    9,
    6,
    8,
    8,
    7,
    6,
    7,
    6,
    5,
    7,
    5,
    8,
    6,
    7,
    8,
    7,
    9,
    7,
    9,
    10,
    6,
    8,
    6,
    7,
    7,
    5,
    7,
    6,
    6,
    5,
    6,
    5,
    7,
    6
)

$jValues = @(
    7,
    8,
    7,
    9,
    7,
    8,
    6,
    5,
    6,
    8,
    8,
    6,
    8,
    6,
    8,
    8,
    8,
    10,
    8,
    7,
    8,
    6,
    8,
    7,
    10,
    10,
    9,
    9,
    6,
    8,
    9,
    7,
    8,
    9,
    7,
    6,
    8,
    6,
    6,
    7,
    6,
    8,
    6,
    7,
    8,
    7,
    9,
    8,
    9,
    10,
    6,
    8,
    7,
    7,
    7,
    7,
    7,
    7,
    7,
    5,
    6,
    5,
    7,
    6
)

for ($k = 0; $k -lt $iValues.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$k]
    $ws.Cells.Item($row, 10).Value = $jValues[$k]
}
